$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the shared formula range from C5:C10 to C5:C11 by filling C11 with the
# same relative formula as the rest of the column.
$ws.Range("C11").Formula = "=B11/200*100"

# Add the new data row
$ws.Range("A11").Value = 45794
$ws.Range("A11").NumberFormat = "d-mmm"
$ws.Range("B11").Value = 62

# Update selection to match the diff (active cell C11)
$ws.Range("C11").Select()
